$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-19 20:34:52"
$wsZh.Range("H4").Value = "2016-03-19 20:35:14"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-19 20:34:55"
$wsDe.Range("H4").Value = "2016-03-19 20:35:19"
